# Update "Table 1" sheet: relabel header, rewrite values as percentage
# strings, add a new bordered (blank) cell at A16, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1")

# Header row
$ws.Range("F1").Value = "Dependent in Last Year"

# Alcohol
$ws.Range("B2").Value = "37.9%"
$ws.Range("C2").Value = "59.2%"
$ws.Range("D2").Value = "13.3%"
$ws.Range("E2").Value = "30.2%"
$ws.Range("F2").Value = "4.5%"

# Marijuana
$ws.Range("B3").Value = "16.6%"
$ws.Range("C3").Value = "38.0%"
$ws.Range("D3").Value = "6.4%"
$ws.Range("E3").Value = "21.2%"
$ws.Range("F3").Value = "2.1%"

# Cigarettes
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "20.0%"
$ws.Range("D4").Value = "0.4%"
$ws.Range("E4").Value = "8.6%"
$ws.Range("F4").Value = "1.2% (All tobacco)"

# Cocaine
$ws.Range("B5").Value = "0.9%"
$ws.Range("C5").Value = "5.6%"
$ws.Range("D5").Value = "0.0%"
$ws.Range("E5").Value = "-"
$ws.Range("F5").Value = "0.0%"

# Heroin
$ws.Range("B6").Value = "0.0%"
$ws.Range("C6").Value = "2.0%"
$ws.Range("D6").Value = "0.0%"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "0.0%"

# New blank cell with a left border, below the footnotes
$ws.Range("A16").Borders.Item(7).LineStyle = 1

# Move the active selection
[void]$ws.Range("F2").Select()
